$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: the URL shown/linked changes to a new gstatic thumbnail link
$ws.Range("A2").Value = "https://encrypted-tbn0.gstatic.com/images?q=tbn:ANd9GcQ2DojOyfbjvs8IDSMwA3TnZZuPAheTa22qkw&s"

# A4: a brand-new row/cell with its own hyperlink to a ferra.ru image
$ws.Range("A4").Value = "https://www.ferra.ru/imgs/2022/02/01/15/5208913/95814d9ab375488468e6df2d5d74b98be7af0bee.webp"
$ws.Hyperlinks.Add($ws.Range("A4"), "https://www.ferra.ru/imgs/2022/02/01/15/5208913/95814d9ab375488468e6df2d5d74b98be7af0bee.webp")
# Reapply the same "Hyperlink" cell style used by the other link cells (A1:A3)
$ws.Range("A4").Style = $ws.Range("A1").Style

# Update the selected/active cell in the sheet view
$ws.Range("M10").Select()
